$d = $word.ActiveDocument

# --- Change 1: Funding statement paragraph ---
$found1 = $d.Content.Find.Execute(
    "J.M. is funded by the Economic and Social Research Council [ES/K006460/1].",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "JM is funded by the Economic and Social Research Council [ES/K006460/1]. No other author had any specific funding. The funder had no input into any aspect of the paper, including data collection, analysis or interpretation of results.",
    2)
Write-Host "Replace1 found:" $found1

# Delete the now-empty paragraph that immediately follows the funding statement
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "$([char]13)" -and $p.Range.Start -gt 300 -and $p.Range.Start -lt 400) {
        Write-Host "Deleting paragraph" $i "at" $p.Range.Start "-" $p.Range.End
        $p.Range.Delete()
        break
    }
}

# --- Change 2: Covering letter paragraph ---
$found2 = $d.Content.Find.Execute(
    "readers. Unlike existing research which is emerging on this area, our manuscript",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "readers. Our manuscript",
    2)
Write-Host "Replace2 found:" $found2

# Locate "readers. O" then insert a collapsed bookmark right after it (before "ur")
$found3 = $d.Content.Find.Execute("readers. O", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Find3 found:" $found3
if ($found3) {
    $pos = $d.Content.Find.Parent.End
    Write-Host "Position after 'readers. O':" $pos
    $bmRange = $d.Range($pos, $pos)
    $bm = $d.Bookmarks.Add("_GoBack", $bmRange)
    Write-Host "Bookmark added"
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($i -ge 10 -and $i -le 13) {
        Write-Host $i "[" $p.Range.Start "," $p.Range.End "] :" $p.Range.Text
    }
}
